$d = $word.ActiveDocument

$d.Content.Find.Execute("83÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "49÷7=", 2) | Out-Null
$d.Content.Find.Execute("24÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "92÷4=", 2) | Out-Null
$d.Content.Find.Execute("67÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "84÷9=", 2) | Out-Null
$d.Content.Find.Execute("21÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷9=", 2) | Out-Null
$d.Content.Find.Execute("42÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "10÷5=", 2) | Out-Null
$d.Content.Find.Execute("54÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "30÷5=", 2) | Out-Null
$d.Content.Find.Execute("52÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "73÷8=", 2) | Out-Null
$d.Content.Find.Execute("87÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "68÷8=", 2) | Out-Null
$d.Content.Find.Execute("72÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "14÷9=", 2) | Out-Null
$d.Content.Find.Execute("44÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "55÷7=", 2) | Out-Null
$d.Content.Find.Execute("20÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "70÷2=", 2) | Out-Null
$d.Content.Find.Execute("36÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "69÷4=", 2) | Out-Null
$d.Content.Find.Execute("32÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "57÷3=", 2) | Out-Null
$d.Content.Find.Execute("49÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "36÷6=", 2) | Out-Null
$d.Content.Find.Execute("31÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "92÷7=", 2) | Out-Null
$d.Content.Find.Execute("65÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "85÷6=", 2) | Out-Null
$d.Content.Find.Execute("18÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "78÷8=", 2) | Out-Null
$d.Content.Find.Execute("21÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "62÷2=", 2) | Out-Null
$d.Content.Find.Execute("36÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "18÷3=", 2) | Out-Null
$d.Content.Find.Execute("17÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷5=", 2) | Out-Null
$d.Content.Find.Execute("52÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "50÷4=", 2) | Out-Null
$d.Content.Find.Execute("82÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷5=", 2) | Out-Null
$d.Content.Find.Execute("47÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "20÷2=", 2) | Out-Null
$d.Content.Find.Execute("11÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "18÷3=", 2) | Out-Null
$d.Content.Find.Execute("38÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "69÷3=", 2) | Out-Null
